$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quarter headers
$ws.Range("E1").Value = "Q4 '25"
$ws.Range("F1").Value = "Q1 '26"

# New data values for the two new quarter columns
$ws.Range("E3").Value = 1365.8
$ws.Range("F3").Value = 1440.6

$ws.Range("E4").Value = 171.4
$ws.Range("F4").Value = 177.5

$ws.Range("E5").Value = 105.8
$ws.Range("F5").Value = 138.4

$ws.Range("E6").Value = 88.7
$ws.Range("F6").Value = 63.1

$ws.Range("E7").Value = 85.7
$ws.Range("F7").Value = 75.7

# Match number format of existing data columns (0.0)
$ws.Range("E3:F7").NumberFormat = "0.0"

# Update selection to match target state (cursor moved to F8 after data entry)
$null = $ws.Range("F8").Select()
